$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A10, A11, A12: apply the header style (same as A9, style index 4) ---
# Use copy/paste-special(formats) so the existing style slot (s="4") is
# reused instead of Excel creating a brand-new duplicate style entry.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats

# --- C11: text "-3" -> "-1" (must remain a text value, not a number) ---
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "-1"
# Re-apply the original formatting (style index 6) from another cell that
# already carries that exact style while holding text, so no new style
# (e.g. with quotePrefix) gets created.
$ws.Range("D17").Copy()
$ws.Range("C11").PasteSpecial(-4122)  # xlPasteFormats

# --- C12: numeric -36 -> -12 ---
$ws.Range("C12").Value = -12

# --- E12: text "9/140" -> "33/140" ---
$ws.Range("E12").Value = "33/140"

$excel.CutCopyMode = 0
